$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.690.40'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.851.40'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = "'262.53"
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  +0.28%  '
$ws.Range("D7").Value = "'0.5352"
$ws.Range("E7").Value = '  +2.42%  '
$ws.Range("D8").Value = "'0.3187"
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").Value = "'0.06970"
$ws.Range("E9").Value = '  +2.09%  '
$ws.Range("D10").Value = "'18.96"
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("D11").Value = "'0.7731"
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("D12").Value = "'0.07829"
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").Value = '1.860.21'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("D14").Value = "'89.88"
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("D15").Value = "'5.063"
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = "'14.18"
$ws.Range("E16").Value = '  +1.31%  '
$ws.Range("E17").Value = '  +0.32%  '
$ws.Range("D18").Value = "'0.000008009"
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '26.727.79'
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("D21").Value = '2.081.62'
$ws.Range("E21").Value = '  -0.21%  '
$ws.Range("D22").Value = "'4.664"
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = "'6.054"
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D24").Value = "'9.417"
$ws.Range("E24").Value = '  -1.48%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = "'143.20"
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = "'2.214"
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("D27").Value = "'1.698"
$ws.Range("E27").Value = '  +2.26%  '
$ws.Range("D28").Value = "'17.16"
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("D29").Value = "'111.84"
$ws.Range("D30").Value = "'4.335"
$ws.Range("E30").Value = '  +2.80%  '
$ws.Range("D31").Value = "'0.08771"
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = "'4.119"
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("D33").Value = "'0.04883"
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("D34").Value = "'0.7419"
$ws.Range("E34").Value = '  +2.94%  '
$ws.Range("D35").Value = "'1.143"
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").Value = "'2.891"
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("D37").Value = "'3.115"
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").Value = "'2.362"
$ws.Range("E38").Value = '  +6.20%  '
$ws.Range("D39").Value = "'0.01752"
$ws.Range("E39").Value = '  -1.79%  '
$ws.Range("D40").Value = "'0.4843"
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("D41").Value = "'0.9075"
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").Value = "'109.34"
$ws.Range("E42").Value = '  -1.97%  '
$ws.Range("D43").Value = "'5.918"
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = '  +0.28%  '
$ws.Range("D45").Value = "'7.724"
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D46").Value = "'0.4215"
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").Value = "'9.172"
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").Value = "'0.1255"
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").Value = "'35.24"
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("B50").Value = 'EOS'
$ws.Range("C50").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D50").Value = "'0.9028"
$ws.Range("E50").Value = '  +1.00%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.05831"
$ws.Range("E51").Value = '  -2.08%  '
